$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.316.19'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '2.397.66'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '502.76'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.77%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.15'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('E7').Value = '  -0.74%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.558'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').Value = '2.428.53'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0971'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.149'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  -4.73%  '
$ws.Range('D14').Value = '2.834.07'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').Value = '57.182.35'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.84'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000133'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '2.387.47'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.24'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '312.76'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.38'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.995'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '65.31'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.995'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').Value = '2.510.72'
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('E28').Value = '  -6.48%  '
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.48'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.98'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').Value = '0.0₃0732'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.19'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.94%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.14'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.997'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.991'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.13'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.97%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.22'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.40%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.83'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('E41').Value = '  +3.95%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '36.20'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '133.64'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +9.04%  '
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.99'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.60%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '258.54'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0914'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0494'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('E51').Value = '  +0.94%  '
